$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Complete row 9: add the two trailing cells (PriceChange, UpDown)
$ws.Range("X9").Value = -0.38999900000000309
$ws.Range("Y9").Value = "Down"

# Add new row 10 (a new trading day's prediction record)
$ws.Range("A10").NumberFormat = "m/d/yy h:mm"
$ws.Range("A10").Value = 42653.881273148145
$ws.Range("B10").Value = -8
$ws.Range("C10").Value = "Sell"
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = 0
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = 0
$ws.Range("N10").Value = 0
$ws.Range("O10").Value = 0
$ws.Range("P10").Value = "Random"
$ws.Range("Q10").Value = 29.009773492518704
$ws.Range("R10").Value = 0.84
$ws.Range("S10").NumberFormat = "0.00%"
$ws.Range("S10").Value = -0.0136
$ws.Range("T10").NumberFormat = "0.00%"
$ws.Range("T10").Value = -0.03
$ws.Range("U10").Value = 14.53
$ws.Range("V10").Value = "N/A"
$ws.Range("W10").Value = -2
